# Clean up code and fix output
# Adds a new "Yearly demand" worksheet at the end of the workbook, matching
# the structure/style used by the other hourly-data sheets in this workbook
# (header row of hours 0-23 in B1:Y1, a row-index column in A2:A4, and the
# corresponding yearly-demand data values).

$wb = $excel.ActiveWorkbook

# Add the new worksheet and move it to the very end of the tab order.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Yearly demand"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

$ws = $wb.Worksheets.Item("Yearly demand")

# Match the outline / page-setup conventions used by the other sheets in
# this workbook (summary rows below, summary columns to the right).
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Match the standard page margins (in points: 1 inch = 72 points) used by
# the other sheets in this workbook.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row: hours 0 .. 23 across B1:Y1
for ($i = 0; $i -le 23; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $i
}

# Row index column: A2=0, A3=1, A4=2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2

# Data values for each of the 3 rows (B..Y => hours 0..23)
$row2 = @(-32.5,-19.5,-13,-13,-13,142.5,291.5,327,388.5,502,596,670.5,745,651,576.5,502,320.5,139,32,-117,-97.5,-78,-52,-39)
$row3 = @(-32.5,-19.5,-13,0,0,-19.5,0,324,486,648,729,751.5,583,567,333.5,340,243,57.99999999999999,-130,0,0,-78,0,-39)
$row4 = @(-32.5,-19.5,0,0,0,-19.5,0,0,81,324,567,589.5,648,567,324,162,81,0,-130,0,0,0,0,-39)

for ($i = 0; $i -le 23; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
    $ws.Cells.Item(4, $i + 2).Value = $row4[$i]
}

# Apply the same formatting (bold, centered, top-aligned, thin-bordered)
# already used for the header row / index column on the other sheets, by
# copying the format from an existing styled cell rather than re-building
# a brand-new style.
$styleSource = $wb.Worksheets.Item("DG Dispatch").Range("B1")
$styleSource.Copy()
$ws.Range("B1:Y1").PasteSpecial(-4122)
$ws.Range("A2:A4").PasteSpecial(-4122)

$excel.CutCopyMode = $false
